# Scheduled price-data refresh for the Chocobo_Profits workbook.
# Updates currentAveragePrice / LevePrice* / LeveProfit* columns (H:N)
# for items whose market data changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 211.5
$ws.Range("I33").Value = 179.44444
$ws.Range("K33").Value = 179.44444
$ws.Range("M33").Value = 49.55556000000001
$ws.Range("H43").Value = 1871.3846
$ws.Range("J43").Value = 2418.25
$ws.Range("L43").Value = 2418.25
$ws.Range("N43").Value = -2556.25
$ws.Range("H106").Value = 4000
$ws.Range("I106").Value = 4000
$ws.Range("K106").Value = 4000
$ws.Range("M106").Value = -3369
$ws.Range("H107").Value = 1702.2609
$ws.Range("I107").Value = 1616
$ws.Range("J107").Value = 1836.4445
$ws.Range("K107").Value = 1616
$ws.Range("L107").Value = 1836.4445
$ws.Range("M107").Value = 304
$ws.Range("N107").Value = -5676.4445
$ws.Range("H112").Value = 1339
$ws.Range("I112").Value = 747.1429000000001
$ws.Range("J112").Value = 1433.159
$ws.Range("K112").Value = 2241.4287
$ws.Range("L112").Value = 4299.477000000001
$ws.Range("M112").Value = -1133.4287
$ws.Range("N112").Value = -6515.477000000001
$ws.Range("H138").Value = 2319.96
$ws.Range("I138").Value = 1060.1875
$ws.Range("J138").Value = 2559.9167
$ws.Range("K138").Value = 3180.5625
$ws.Range("L138").Value = 7679.750100000001
$ws.Range("M138").Value = 1959.4375
$ws.Range("N138").Value = -17959.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 619.9459000000001
$ws.Range("I2").Value = 567.6667
$ws.Range("J2").Value = 761.1
$ws.Range("K2").Value = 567.6667
$ws.Range("L2").Value = 761.1
$ws.Range("M2").Value = -454.6667
$ws.Range("N2").Value = -987.1
$ws.Range("H32").Value = 5243
$ws.Range("I32").Value = 4475.846
$ws.Range("J32").Value = 8567.333000000001
$ws.Range("K32").Value = 4475.846
$ws.Range("L32").Value = 8567.333000000001
$ws.Range("M32").Value = -4188.846
$ws.Range("N32").Value = -9141.333000000001
$ws.Range("H61").Value = 903.4194
$ws.Range("I61").Value = 903.4194
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 903.4194
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -691.4194
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2529.6086
$ws.Range("I74").Value = 2334.372
$ws.Range("J74").Value = 5328
$ws.Range("K74").Value = 2334.372
$ws.Range("L74").Value = 5328
$ws.Range("M74").Value = -1460.372
$ws.Range("N74").Value = -7076
$ws.Range("H77").Value = 2529.6086
$ws.Range("I77").Value = 2334.372
$ws.Range("J77").Value = 5328
$ws.Range("K77").Value = 11671.86
$ws.Range("L77").Value = 26640
$ws.Range("M77").Value = -7303.859999999999
$ws.Range("N77").Value = -35376
$ws.Range("H104").Value = 34500
$ws.Range("J104").Value = 34500
$ws.Range("L104").Value = 34500
$ws.Range("N104").Value = -41488
$ws.Range("H116").Value = 619.9459000000001
$ws.Range("I116").Value = 567.6667
$ws.Range("J116").Value = 761.1
$ws.Range("K116").Value = 567.6667
$ws.Range("L116").Value = 761.1
$ws.Range("M116").Value = 1726.3333
$ws.Range("N116").Value = -5349.1
$ws.Range("H132").Value = 2333.4814
$ws.Range("I132").Value = 1295.5238
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 3886.5714
$ws.Range("L132").Value = 17899.0005
$ws.Range("M132").Value = -1356.5714
$ws.Range("N132").Value = -22959.0005
$ws.Range("H133").Value = 39260
$ws.Range("J133").Value = 39260
$ws.Range("L133").Value = 39260
$ws.Range("N133").Value = -44320
$ws.Range("H136").Value = 903.4194
$ws.Range("I136").Value = 903.4194
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2710.2582
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -160.2582000000002
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 619.9459000000001
$ws.Range("I3").Value = 567.6667
$ws.Range("J3").Value = 761.1
$ws.Range("K3").Value = 567.6667
$ws.Range("L3").Value = 761.1
$ws.Range("M3").Value = -453.6667
$ws.Range("N3").Value = -989.1
$ws.Range("H99").Value = 1877.0588
$ws.Range("I99").Value = 1064.5454
$ws.Range("K99").Value = 1064.5454
$ws.Range("M99").Value = 433.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6945638
$ws.Range("I16").Value = 11112122
$ws.Range("K16").Value = 11112122
$ws.Range("M16").Value = -11111835
$ws.Range("H31").Value = 10418335
$ws.Range("I31").Value = 789.13513
$ws.Range("J31").Value = 45459172
$ws.Range("K31").Value = 789.13513
$ws.Range("L31").Value = 45459172
$ws.Range("M31").Value = -494.13513
$ws.Range("N31").Value = -45459762
$ws.Range("H34").Value = 10418335
$ws.Range("I34").Value = 789.13513
$ws.Range("J34").Value = 45459172
$ws.Range("K34").Value = 789.13513
$ws.Range("L34").Value = 45459172
$ws.Range("M34").Value = -587.13513
$ws.Range("N34").Value = -45459576
$ws.Range("H58").Value = 1779.421
$ws.Range("I58").Value = 1621.5646
$ws.Range("J58").Value = 2478.5
$ws.Range("K58").Value = 1621.5646
$ws.Range("L58").Value = 2478.5
$ws.Range("M58").Value = -1418.5646
$ws.Range("N58").Value = -2884.5
$ws.Range("H86").Value = 2555.375
$ws.Range("I86").Value = 1528.8
$ws.Range("J86").Value = 4266.3335
$ws.Range("K86").Value = 1528.8
$ws.Range("L86").Value = 4266.3335
$ws.Range("M86").Value = -405.8
$ws.Range("N86").Value = -6512.3335
$ws.Range("H89").Value = 2555.375
$ws.Range("I89").Value = 1528.8
$ws.Range("J89").Value = 4266.3335
$ws.Range("K89").Value = 7644
$ws.Range("L89").Value = 21331.6675
$ws.Range("M89").Value = -2028
$ws.Range("N89").Value = -32563.6675
$ws.Range("H113").Value = 6945638
$ws.Range("I113").Value = 11112122
$ws.Range("K113").Value = 11112122
$ws.Range("M113").Value = -11109952
$ws.Range("H136").Value = 1779.421
$ws.Range("I136").Value = 1621.5646
$ws.Range("J136").Value = 2478.5
$ws.Range("K136").Value = 4864.6938
$ws.Range("L136").Value = 7435.5
$ws.Range("M136").Value = -2314.6938
$ws.Range("N136").Value = -12535.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 467
$ws.Range("I92").Value = 451
$ws.Range("J92").Value = 499
$ws.Range("K92").Value = 1353
$ws.Range("L92").Value = 1497
$ws.Range("M92").Value = -105
$ws.Range("N92").Value = -3993
$ws.Range("H113").Value = 583.6429000000001
$ws.Range("I113").Value = 604.13336
$ws.Range("J113").Value = 560
$ws.Range("K113").Value = 1812.40008
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = 357.5999199999999
$ws.Range("N113").Value = -6020
$ws.Range("H131").Value = 9260178
$ws.Range("J131").Value = 958.96075
$ws.Range("L131").Value = 2876.88225
$ws.Range("N131").Value = -12956.88225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2309.72
$ws.Range("I102").Value = 1770.1364
$ws.Range("J102").Value = 6266.6665
$ws.Range("K102").Value = 1770.1364
$ws.Range("L102").Value = 6266.6665
$ws.Range("M102").Value = -148.1364000000001
$ws.Range("N102").Value = -9510.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1539.174
$ws.Range("I61").Value = 1454.6
$ws.Range("J61").Value = 1697.75
$ws.Range("K61").Value = 1454.6
$ws.Range("L61").Value = 1697.75
$ws.Range("M61").Value = -1252.6
$ws.Range("N61").Value = -2101.75
$ws.Range("H113").Value = 1539.174
$ws.Range("I113").Value = 1454.6
$ws.Range("J113").Value = 1697.75
$ws.Range("K113").Value = 1454.6
$ws.Range("L113").Value = 1697.75
$ws.Range("M113").Value = 715.4000000000001
$ws.Range("N113").Value = -6037.75
